$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K" = strikeouts),
# replacing the previous "Strike#" counts per the regenerated save_data.
$gValues = @{
    2 = 0;
    3 = 2;
    4 = 0;
    5 = 2;
    6 = 1;
    7 = 1;
    8 = 0;
    9 = 0;
    10 = 1;
    11 = 1;
    12 = 0;
    13 = 2;
    14 = 2;
    15 = 2;
    16 = 5;
    18 = 1;
    19 = 0;
    20 = 5;
    21 = 1;
    22 = 3;
    23 = 0;
    24 = 2;
    25 = 1;
    26 = 1;
    27 = 2;
    28 = 0;
    29 = 2;
    30 = 3;
    31 = 2;
    32 = 1;
    33 = 2;
    34 = 2;
    35 = 1;
    36 = 1;
    37 = 0;
    38 = 2;
    39 = 1;
    40 = 2;
    41 = 1;
    42 = 2;
    43 = 1;
    44 = 0;
    45 = 2;
    46 = 2;
    47 = 1;
    48 = 3;
    49 = 2;
    50 = 1;
    51 = 1;
    52 = 0;
    53 = 1;
    54 = 0;
    55 = 0;
    56 = 3;
    57 = 1;
    58 = 1;
    59 = 1;
    60 = 1;
    61 = 0;
    62 = 1;
    63 = 2;
    64 = 1;
    65 = 1;
    66 = 2;
    68 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
